{"js": "// \"Add design plantilla docente V2\"\n// 1) Resize/reposition the masthead logo (legacy VML v:shape inside a w:pict).\n// 2) Tighten the leading spaces on the two title lines.\n// 3) Re-space a handful of the \"ficha\" fields (\u00c1rea/LOCALIDAD, CORREO/ACUERDO,\n//    REPRESENTANTE LEGAL/CELULAR).\n\n// --- 1) VML shape (logo) style change -------------------------------------\n// The picture is a legacy VML <v:shape> wrapped in <w:pict>; it is not\n// surfaced through body.inlinePictures/shapes, so we go through the raw\n// paragraph OOXML (getOoxml/insertOoxml) and only touch the style string -\n// everything else (the image relationship, etc.) is left untouched.\nconst shapeParagraph = context.document.body.paragraphs.getFirst();\nconst shapeOoxmlResult = shapeParagraph.getOoxml();\nawait context.sync();\n\nconst oldShapeStyle =\n  \"width:130pt; height:65pt; margin-left:-1pt; margin-top:-1pt; mso-position-horizontal:left; mso-position-vertical:top; mso-position-horizontal-relative:char; mso-position-vertical-relative:line; z-index:-2147483647;\";\nconst newShapeStyle =\n  \"width:165pt; height:75pt; margin-left:90pt; margin-top:50pt; position:absolute; mso-position-horizontal:left; mso-position-vertical:top; mso-position-horizontal-relative:char; mso-position-vertical-relative:line; z-index:-2147483647;\";\n\nlet shapeXml = shapeOoxmlResult.value;\nif (shapeXml.indexOf(oldShapeStyle) === -1) {\n  throw new Error(\"Could not locate the expected v:shape style to update.\");\n}\nshapeXml = shapeXml.replace(oldShapeStyle, newShapeStyle);\nshapeParagraph.insertOoxml(shapeXml, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2) & 3) plain text re-spacing -----------------------------------------\nconst textReplacements = [\n  {\n    old:\n      \"                                                           DIRRECCI\u00d3N GENERAL DE EDUCACI\u00d3N B\u00c1SICA\",\n    new:\n      \"                                                        DIRRECCI\u00d3N GENERAL DE EDUCACI\u00d3N B\u00c1SICA\",\n  },\n  {\n    old:\n      \"                                                           DEPARTAMENTO DE ESCUELAS PARTICULARES\",\n    new:\n      \"                                                        DEPARTAMENTO DE ESCUELAS PARTICULARES\",\n  },\n  {\n    old:\n      \"Area: 02 LOCALIDAD:HUNUCMA                                                                    MUNICIPIO:HUNUCMA                               TELEFONO DE CT: 98893110737\",\n    new:\n      \"\u00c1rea: 02   LOCALIDAD: HUNUCMA                                                                    MUNICIPIO:HUNUCMA                               TELEFONO DE CT: 98893110737\",\n  },\n  {\n    old:\n      \"CORREO ELECTRONICO DEL CT: fray_diego69@hotmail.com                                   NO\u00b0 ACUERDO 208         FECHA ACUERDO 29/07/199\",\n    new:\n      \"CORREO ELECTRONICO DEL CT: fray_diego69@hotmail.com                 NO\u00b0 ACUERDO 208                                              FECHA ACUERDO 29/07/199\",\n  },\n  {\n    old:\n      \"REPRESENTANTE LEGAL: R\u0332O\u0332L\u0332A\u0332N\u0332D\u0332O\u0332_\u0332J\u0332A\u0332V\u0332I\u0332E\u0332R\u0332 Q\u0332U\u0332I\u0332N\u0332T\u0332A\u0332L\u0332_\u0332C\u0332A\u0332S\u0332T\u0332I\u0332L\u0332L\u0332A\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_                      CELULAR DEL REPRESENTANTE:9999000667\",\n    new:\n      \"REPRESENTANTE LEGAL: R\u0332O\u0332L\u0332A\u0332N\u0332D\u0332O\u0332_\u0332J\u0332A\u0332V\u0332I\u0332E\u0332R\u0332 Q\u0332U\u0332I\u0332N\u0332T\u0332A\u0332L\u0332_\u0332C\u0332A\u0332S\u0332T\u0332I\u0332L\u0332L\u0332A\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_                                 CELULAR DEL REPRESENTANTE:9999000667\",\n  },\n];\n\nfor (const { old: oldText, new: newText } of textReplacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text to replace not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Masthead logo (legacy VML v:shape inside a w:pict): resize/reposition\n#    by rewriting just the style attribute via the paragraph's WordOpenXML.\n$oldShapeTag = '<v:shape type=\"#_x0000_t75\" style=\"width:130pt; height:65pt; margin-left:-1pt; margin-top:-1pt; mso-position-horizontal:left; mso-position-vertical:top; mso-position-horizontal-relative:char; mso-position-vertical-relative:line; z-index:-2147483647;\">'\n$newShapeTag = '<v:shape type=\"#_x0000_t75\" style=\"width:165pt; height:75pt; margin-left:90pt; margin-top:50pt; position:absolute; mso-position-horizontal:left; mso-position-vertical:top; mso-position-horizontal-relative:char; mso-position-vertical-relative:line; z-index:-2147483647;\">'\n$shapePara = $d.Paragraphs(1).Range\n$shapeXml = $shapePara.WordOpenXML\nif ($shapeXml -notlike \"*$oldShapeTag*\") { throw 'Expected v:shape style not found' }\n$shapeXml = $shapeXml.Replace($oldShapeTag, $newShapeTag)\n$shapePara.InsertXML($shapeXml)\n\n# 2) Plain-text re-spacing via Find/Replace (exact, non-wildcard).\n$old = '                                                           DIRRECCI\u00d3N GENERAL DE EDUCACI\u00d3N B\u00c1SICA'\n$new = '                                                        DIRRECCI\u00d3N GENERAL DE EDUCACI\u00d3N B\u00c1SICA'\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $old\n$find.Replacement.Text = $new\n$found = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\nif (-not $found) { throw \"Text not found: $old\" }\n\n$old = '                                                           DEPARTAMENTO DE ESCUELAS PARTICULARES'\n$new = '                                                        DEPARTAMENTO DE ESCUELAS PARTICULARES'\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $old\n$find.Replacement.Text = $new\n$found = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\nif (-not $found) { throw \"Text not found: $old\" }\n\n$old = 'Area: 02 LOCALIDAD:HUNUCMA                                                                    MUNICIPIO:HUNUCMA                               TELEFONO DE CT: 98893110737'\n$new = '\u00c1rea: 02   LOCALIDAD: HUNUCMA                                                                    MUNICIPIO:HUNUCMA                               TELEFONO DE CT: 98893110737'\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $old\n$find.Replacement.Text = $new\n$found = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\nif (-not $found) { throw \"Text not found: $old\" }\n\n$old = 'CORREO ELECTRONICO DEL CT: fray_diego69@hotmail.com                                   NO\u00b0 ACUERDO 208         FECHA ACUERDO 29/07/199'\n$new = 'CORREO ELECTRONICO DEL CT: fray_diego69@hotmail.com                 NO\u00b0 ACUERDO 208                                              FECHA ACUERDO 29/07/199'\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $old\n$find.Replacement.Text = $new\n$found = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\nif (-not $found) { throw \"Text not found: $old\" }\n\n$old = 'REPRESENTANTE LEGAL: R\u0332O\u0332L\u0332A\u0332N\u0332D\u0332O\u0332_\u0332J\u0332A\u0332V\u0332I\u0332E\u0332R\u0332 Q\u0332U\u0332I\u0332N\u0332T\u0332A\u0332L\u0332_\u0332C\u0332A\u0332S\u0332T\u0332I\u0332L\u0332L\u0332A\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_                      CELULAR DEL REPRESENTANTE:9999000667'\n$new = 'REPRESENTANTE LEGAL: R\u0332O\u0332L\u0332A\u0332N\u0332D\u0332O\u0332_\u0332J\u0332A\u0332V\u0332I\u0332E\u0332R\u0332 Q\u0332U\u0332I\u0332N\u0332T\u0332A\u0332L\u0332_\u0332C\u0332A\u0332S\u0332T\u0332I\u0332L\u0332L\u0332A\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_\u0332_                                 CELULAR DEL REPRESENTANTE:9999000667'\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $old\n$find.Replacement.Text = $new\n$found = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\nif (-not $found) { throw \"Text not found: $old\" }\n\n"}
